$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 59 - new entry: 4 botellones
$ws.Range("B59").Value = 45254
$ws.Range("C59").Value = "4 botellones"
$ws.Range("D59").Value = -212

# Row 60 - new entry: 3 botellones
$ws.Range("B60").Value = 45258
$ws.Range("C60").Value = "3 botellones"
$ws.Range("D60").Value = -159

# Update selection to reflect next empty entry row
$ws.Range("B61").Select()
